# Q3 Update - 2025
# Applies the UN-BOT.xlsx Q3 refresh:
#   1. The global "short-url" token in column B changes for every data row.
#   2. The last "year" block (coo = Botswana group, rows 443-453) gets refreshed
#      figures, and the dataset shrinks from 453 to 450 data rows (the old
#      rows 451-453 disappear entirely once the new, shorter block lands).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a value into a cell while preserving its original "Text"
# storage type. The source file stores every cell (numbers included) as
# shared-string text; Excel's default Value setter auto-converts
# numeric-looking strings to real numbers, so for anything that looks like a
# number we force the cell to Text format first.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($Sheet, [string]$Addr, [string]$Val)

    $looksNumeric = $Val -match '^-?\d+(\.\d+)?$'
    if ($looksNumeric) {
        $Sheet.Range($Addr).NumberFormat = "@"
    }
    $Sheet.Range($Addr).Value = $Val
}

# ---------------------------------------------------------------------------
# 1) Remove the trailing rows 451-453 (the final 3 rows of the "Botswana"
#    block are dropped in the refreshed dataset). Deleting first so the row
#    numbers used below (443-450) are unaffected either way.
# ---------------------------------------------------------------------------
$ws.Rows("451:453").Delete()

# ---------------------------------------------------------------------------
# 2) Refresh figures for rows 443-450.
# ---------------------------------------------------------------------------
Set-TextValue $ws "T443" "53"

Set-TextValue $ws "N444" "46"
Set-TextValue $ws "P444" "7"

Set-TextValue $ws "T445" "29"

Set-TextValue $ws "N446" "353"
Set-TextValue $ws "O446" "7"

Set-TextValue $ws "F447" "161"
Set-TextValue $ws "G447" "Rwanda"
Set-TextValue $ws "H447" "RWA"
Set-TextValue $ws "I447" "RWA"
Set-TextValue $ws "N447" "18"

Set-TextValue $ws "F448" "172"
Set-TextValue $ws "G448" "Somalia"
Set-TextValue $ws "H448" "SOM"
Set-TextValue $ws "I448" "SOM"
Set-TextValue $ws "N448" "375"

Set-TextValue $ws "F449" "199"
Set-TextValue $ws "G449" "Uganda"
Set-TextValue $ws "H449" "UGA"
Set-TextValue $ws "I449" "UGA"
Set-TextValue $ws "N449" "8"
Set-TextValue $ws "O449" "0"

Set-TextValue $ws "F450" "214"
Set-TextValue $ws "G450" "Zimbabwe"
Set-TextValue $ws "H450" "ZIM"
Set-TextValue $ws "I450" "ZWE"
Set-TextValue $ws "N450" "16"
Set-TextValue $ws "T450" "37"

# ---------------------------------------------------------------------------
# 3) Global short-url refresh: every data row (2-450 after the deletion
#    above) gets the new short-url token in column B.
# ---------------------------------------------------------------------------
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "n8hFl7"
}

Write-Output "Q3 2025 update applied"
